$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.339.75'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '1.823.45'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.10'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4464'
$ws.Range("E7").Value = '  -2.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3751'
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07469'
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("E10").Value = '  +2.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.97'
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").Value = '1.831.71'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.756'
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.421'
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.64'
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07128'
$ws.Range("E16").Value = '  +0.79%  '
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008764'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.14'
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").Value = '27.342.71'
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.418'
$ws.Range("E22").Value = '  +4.20%  '
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("D24").Value = '2.057.29'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.959'
$ws.Range("E25").Value = '  -1.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.31'
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.310'
$ws.Range("E27").Value = '  +2.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.65'
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.378'
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.85'
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08888'
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7851'
$ws.Range("E32").Value = '  +3.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.204'
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.617'
$ws.Range("E34").Value = '  +3.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.910'
$ws.Range("E35").Value = '  -2.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9996'
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.110'
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01992'
$ws.Range("E38").Value = '  +1.02%  '
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.305'
$ws.Range("E40").Value = '  +1.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5351'
$ws.Range("E41").Value = '  -0.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.855'
$ws.Range("E42").Value = '  -1.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1715'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.291'
$ws.Range("E44").Value = '  +16.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.658'
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5107'
$ws.Range("E46").Value = '  -2.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.59'
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.692'
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.11'
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9995'
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06394'
$ws.Range("E51").Value = '  +0.56%  '
